$d = $word.ActiveDocument

function Set-ParagraphRuns($paragraphIndex, $styleId, $words) {
    $runsXml = ""
    foreach ($w in $words) {
        $runsXml += "<w:r><w:t xml:space=`"preserve`">$w</w:t></w:r>"
    }

    $xml = '<?xml version="1.0" standalone="yes"?>' +
        '<?mso-application progid="Word.Document"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:body>' +
        '<w:p><w:pPr><w:pStyle w:val="' + $styleId + '"/></w:pPr>' + $runsXml + '</w:p>' +
        '</w:body></w:document>' +
        '</pkg:xmlData></pkg:part></pkg:package>'

    $range = $d.Paragraphs($paragraphIndex).Range
    $range.InsertXML($xml)
}

# Title: "Answers: Rationalizing the denominator"
Set-ParagraphRuns 1 "Title" @("Answers:", " ", "Rationalizing", " ", "the", " ", "denominator")

# Author: "Maximilian Volmar"
Set-ParagraphRuns 2 "Author" @("Maximilian", " ", "Volmar")

# Abstract: "Answers to questions relating to the guide on rationalizing the denominator."
Set-ParagraphRuns 4 "Abstract" @("Answers", " ", "to", " ", "questions", " ", "relating", " ", "to", " ", "the", " ", "guide", " ", "on", " ", "rationalizing", " ", "the", " ", "denominator.")
